# Generate Report for Archive
#
# The localization status for this item moved on from "Ready for handoff"
# to "In Translation" everywhere it is reported: the per-language status
# cells on the "Overview" sheet (E2 for zh-cn, F2 for de-de) and the
# "Status" column on each language detail sheet (C2 on "zh-cn" and
# "de-de"). Shortening that text also lets Excel's status column shrink
# to fit the new, shorter label.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Target column width (characters) for the narrower "Status" columns,
# taken straight from the saved workbook after the text shrank. Excel's
# ColumnWidth setter snaps to whole on-screen pixels (width*6 + 5px
# padding), so back the literal target out of that same 5px/6 offset.
$targetColumnWidth = 13.4101845877511 - (5 / 6)

# --- Overview sheet: per-language status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth

# --- zh-cn detail sheet: Status column ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth

# --- de-de detail sheet: Status column ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
